$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B1/C1 header labels (LastName <-> FirstName) and clear their
# explicit cell style so they fall back to the default "Normal" style.
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("B1:C1").Style = "Normal"

# Overwrite the CNE identifiers in column A (rows 2-11) with the corrected
# value.
$ws.Range("A2:A11").Value = 18000031

# Move the active selection to F9.
$ws.Range("F9").Select()
